$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shared string needed: "ECs" (used for Target cluster in new rows)

# Row 2: Target cluster FAPs -> ECs, and updated values
$ws.Range("A2:T7").Value = $null

$rows = @(
  @("FAPs","Fgf16","Fgfr2","ECs",3,1,1.678104,5.034312,0.551436927751233,0.551436927751233,1,0.3333333333333333,0.4871643333333333,1.461493,0.1914458468300136,0.1914458468300136,0.817512416424,7.357611747816,0.1055703096066758,0.1055703096066758),
  @("FAPs","Fgf16","Fgfr2","FAPs",3,1,1.678104,5.034312,0.551436927751233,0.551436927751233,3,1,1.864751,5.594253,0.7328098752209857,0.7328098752209857,3.129246112104,28.163215008936,0.4040984262176248,0.4040984262176248),
  @("FAPs","Fgf16","Fgfr2","sCs",3,1,1.678104,5.034312,0.551436927751233,0.551436927751233,3,1,0.1927433333333333,0.57823,0.07574427794900063,0.07574427794900063,0.32344335864,2.91099022776,0.04176819192693237,0.04176819192693237),
  @("sCs","Fgf16","Fgfr2","ECs",3,1,1.365043666666667,4.095131,0.448563072248767,0.448563072248767,1,0.3333333333333333,0.4871643333333333,1.461493,0.1914458468300136,0.1914458468300136,0.6650005878425556,5.985005290583,0.08587553722333775,0.08587553722333777),
  @("sCs","Fgf16","Fgfr2","FAPs",3,1,1.365043666666667,4.095131,0.448563072248767,0.448563072248767,3,1,1.864751,5.594253,0.7328098752209857,0.7328098752209857,2.545466542460334,22.909198882143,0.3287114490033609,0.328711449003361),
  @("sCs","Fgf16","Fgfr2","sCs",3,1,1.365043666666667,4.095131,0.448563072248767,0.448563072248767,3,1,0.1927433333333333,0.57823,0.07574427794900063,0.07574427794900063,0.2631030664588889,2.36792759813,0.03397608602206825,0.03397608602206826)
)

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
